# PLAU pagination xpath fix
# Adds a new "PLNewAUtestuser4" row to the Users sheet (row 42), mirroring
# the existing PLAUtestuser4 row (row 41), with hyperlinks on the email
# columns (G and K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$newRow = 42

# --- Column A: UserName ---
$ws.Cells.Item($newRow, 1).Value = "PLNewAUtestuser4"

# --- Column B: Password ---
$ws.Cells.Item($newRow, 2).Value = "P@ssword2"

# --- Column G: Email (hyperlink) ---
$email = "PLNewAUtestuser4@mailinator.com"
$ws.Cells.Item($newRow, 7).Value = $email
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 7), "mailto:$email") | Out-Null

# --- Column J: ProdUserName ---
$ws.Cells.Item($newRow, 10).Value = "ProdAUtestuser4"

# --- Column K: ProdEmail (hyperlink) ---
$prodEmail = "ProdAUtestuser4@mailinator.com"
$ws.Cells.Item($newRow, 11).Value = $prodEmail
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 11), "mailto:$prodEmail") | Out-Null

# --- Column L: password hint ---
$ws.Cells.Item($newRow, 12).Value = "thomsonreuters"

# Update the view: scroll down a bit and move the frozen-pane split / selection
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Application.ActiveWindow.SplitColumn = 5
$ws.Range("H50").Select()
